$wb = $excel.ActiveWorkbook

# --- Rename header labels on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the existing sheets ---
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Header row
$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($col = 1; $col -le 4; $col++) {
    $cell = $wsForecast.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows
$data = @(
    @(45354.99999999999, 13, -0.8017563315835705, 27.72524557847934),
    @(45361.99999999999, 13, -0.4898744285013225, 28.51305206404469),
    @(45382.99999999999, 12, -3.258397315292796, 26.557588587565),
    @(45410.99999999999, 11, -2.738414766955998, 25.78200658238674),
    @(45417.99999999999, 11, -3.58060771277047, 25.64528369526806),
    @(45424.99999999999, 10, -4.411352587027586, 24.19946502390025),
    @(45431.99999999999, 10, -5.03609465226006, 24.42730509745734),
    @(45438.99999999999, 10, -4.739800410958991, 23.39248955715412),
    @(45445.99999999999, 9, -5.577046478275693, 23.80874163806497),
    @(45452.99999999999, 9, -5.941389722250211, 23.39324616834318),
    @(45459.99999999999, 9, -5.389238000103989, 22.51328344172301),
    @(45466.99999999999, 9, -6.510894222836977, 22.62552815442562),
    @(45487.99999999999, 8, -6.426886711781551, 23.30543344165861),
    @(45494.99999999999, 7, -6.91033985600642, 21.60333499988614),
    @(45501.99999999999, 7, -7.320239562221132, 21.34477540803483),
    @(45508.99999999999, 7, -7.318542692272477, 20.66148188671337),
    @(45515.99999999999, 6, -7.308447893041321, 21.15945929185906),
    @(45522.99999999999, 6, -7.476973474343501, 19.87957213923579),
    @(45529.99999999999, 6, -9.173493661536215, 19.87727790046439),
    @(45536.99999999999, 5, -9.282956299689527, 20.66982326624042),
    @(45543.99999999999, 5, -8.739187262192258, 19.80680294824181)
)

$r = 2
foreach ($row in $data) {
    $dateCell = $wsForecast.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]

    $r = $r + 1
}
